$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, positioned right before "总计".
#    Copy "2021-Q4" as a template so the new sheet inherits the exact
#    same header/row styling (bold header, thin borders, centered).
# ------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template ("2021-Q4") has 4 data rows; "2022-Q1" only needs 2, so
# drop the extra two rows.
$newSheet.Rows("4:5").Clear()

# Make sure the fund-code / numeric-looking text columns stay TEXT
# (otherwise values like "011052" or "0.20" get coerced to numbers
# and lose their leading zero / trailing zero).
$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "011052"
$newSheet.Range("C2").Value = "鹏华弘裕一年持有期混合A"
$newSheet.Range("D2").Value = "2.92"
$newSheet.Range("E2").Value = "24.56"
$newSheet.Range("F2").Value = "1.45"
$newSheet.Range("G2").Value = "0.0423"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "011053"
$newSheet.Range("C3").Value = "鹏华弘裕一年持有期混合C"
$newSheet.Range("D3").Value = "0.20"
$newSheet.Range("E3").Value = "24.56"
$newSheet.Range("F3").Value = "1.45"
$newSheet.Range("G3").Value = "0.0029"
$newSheet.Range("H3").Value = 5

# The data cells themselves carry no special formatting in the
# original sheets (only the header row / index column do) - drop the
# "@" override now that the values are safely stored as text.
$newSheet.Range("B2:G3").ClearFormats()

# ------------------------------------------------------------------
# 2) Add the 2022-Q1 summary row to "总计" (insert above the existing
#    2021-Q4 row, then renumber the index column).
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows("2:2").Insert()

# Match the formatting of the surrounding data rows exactly (the bare
# Insert() leaves the new row with row-1's header style).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.05

# Renumber the remaining index column (0,1,2,3) now that a row was
# inserted in front of them.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
